$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 24, pushing all existing
# data (old rows 24-94) down to rows 26-96. This mirrors the weekly
# update pattern already present in the sheet (newest record inserted
# at the top of the data block).
$ws.Rows("24:25").Insert()

# Populate the two freshly inserted rows with the new week's records.

# Row 24: "Primera" quality record dated 2021-11-25 (serial 44525)
$ws.Range("A24").Value = 2
$ws.Range("B24").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44525
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 100112043
$ws.Range("G24").Value = "Pepino ensalada"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 6500
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 6750
$ws.Range("N24").Value = "`$/caja 70 unidades"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 96
$ws.Range("Q24").Value = 70
$ws.Range("R24").Value = "Hortaliza"

# Row 25: "Segunda" quality record, same date
$ws.Range("A25").Value = 2
$ws.Range("B25").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 44525
$ws.Range("E25").Value = 4
$ws.Range("F25").Value = 100112043
$ws.Range("G25").Value = "Pepino ensalada"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 4500
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = 4750
$ws.Range("N25").Value = "`$/caja 100 unidades"
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 48
$ws.Range("Q25").Value = 100
$ws.Range("R25").Value = "Hortaliza"
